$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Milestone list. Entered out of row order (B7 last) to reproduce the
# shared-string table ordering of the source workbook.
$ws.Range("B3").Value = "Define functionality"
$ws.Range("B4").Value = "Find core components"
$ws.Range("B5").Value = "Schematics"
$ws.Range("B6").Value = "Layout"
$ws.Range("B8").Value = "Production files generation"
$ws.Range("B9").Value = "Board house order"
$ws.Range("B10").Value = "Arduino board file"
$ws.Range("B7").Value = "Mechanical integration"
$ws.Range("B11").Value = "Programming!"

# Column B was widened to fit the longer milestone text
$ws.Columns("B").ColumnWidth = 26

# Selection left on B19 (matches the saved view state)
$null = $ws.Range("B19").Select()
